$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the "-" placeholder values for the new row
$ws.Range("A2").Value = "-"
$ws.Range("B2").Value = "-"

# Update the active cell/selection to B3
$ws.Range("B3").Select()
